# Macroferia Regional de Talca - Naranja: weekly fruit/vegetable update.
#
# The former rows 255-256 (Fukumoto entries dated 44399) are being
# superseded by new price readings dated 44448 (Lane Late / Navel Late),
# while the original 255-256 rows are preserved by shifting them down to
# rows 257-258. Two brand-new rows (261-262) are appended at the end for
# date 44400.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Make room: insert two blank rows before the current row 257. This
#    pushes the existing rows 257-260 down to 259-262, and leaves rows
#    255-256 (soon to be edited) untouched above the insertion point.
$ws.Rows("257:258").Insert()

# 2) Re-populate the two freshly inserted rows (257-258) with the content
#    that used to live in rows 255-256 before this edit (unchanged data).
$ws.Range("A257").Value = 5
$ws.Range("B257").Value = "Macroferia Regional de Talca"
$ws.Range("C257").Value = "Maule"
$ws.Range("D257").Value = 44399
$ws.Range("E257").Value = 7
$ws.Range("F257").Value = "Fruta"
$ws.Range("G257").Value = 100102
$ws.Range("H257").Value = "Cítricos"
$ws.Range("I257").Value = 100102005
$ws.Range("J257").Value = "Naranja"
$ws.Range("K257").Value = "Fukumoto"
$ws.Range("L257").Value = "Primera"
$ws.Range("M257").Value = 300
$ws.Range("N257").Value = 5500
$ws.Range("O257").Value = 5500
$ws.Range("P257").Value = 5500
$ws.Range("Q257").Value = "$/bandeja 15 kilos granel"
$ws.Range("R257").Value = "Provincia de Melipilla"
$ws.Range("S257").Value = 367
$ws.Range("T257").Value = 15

$ws.Range("A258").Value = 5
$ws.Range("B258").Value = "Macroferia Regional de Talca"
$ws.Range("C258").Value = "Maule"
$ws.Range("D258").Value = 44399
$ws.Range("E258").Value = 7
$ws.Range("F258").Value = "Fruta"
$ws.Range("G258").Value = 100102
$ws.Range("H258").Value = "Cítricos"
$ws.Range("I258").Value = 100102005
$ws.Range("J258").Value = "Naranja"
$ws.Range("K258").Value = "Fukumoto"
$ws.Range("L258").Value = "Primera"
$ws.Range("M258").Value = 500
$ws.Range("N258").Value = 6000
$ws.Range("O258").Value = 6000
$ws.Range("P258").Value = 6000
$ws.Range("Q258").Value = "$/bandeja 15 kilos granel"
$ws.Range("R258").Value = "Región de O'Higgins"
$ws.Range("S258").Value = 400
$ws.Range("T258").Value = 15

# 3) Update the original rows 255-256 with this week's new readings.
$ws.Range("D255").Value = 44448
$ws.Range("K255").Value = "Lane Late"
$ws.Range("N255").Value = 6000
$ws.Range("O255").Value = 6000
$ws.Range("P255").Value = 6000
$ws.Range("Q255").Value = "$/bandeja 18 kilos granel"
$ws.Range("S255").Value = 333
$ws.Range("T255").Value = 18

$ws.Range("D256").Value = 44448
$ws.Range("K256").Value = "Navel Late"
$ws.Range("M256").Value = 440
$ws.Range("N256").Value = 5500
$ws.Range("O256").Value = 6000
$ws.Range("P256").Value = 5705
$ws.Range("S256").Value = 380

# 4) Fix up the (now shifted) rows 259-260, which keep their own data but
#    need the date correction from 44400 back to 44399.
$ws.Range("D259").Value = 44399
$ws.Range("D260").Value = 44399

# 5) Append two brand-new rows (261-262) with fresh data points.
$ws.Range("A261").Value = 5
$ws.Range("B261").Value = "Macroferia Regional de Talca"
$ws.Range("C261").Value = "Maule"
$ws.Range("D261").Value = 44400
$ws.Range("E261").Value = 7
$ws.Range("F261").Value = "Fruta"
$ws.Range("G261").Value = 100102
$ws.Range("H261").Value = "Cítricos"
$ws.Range("I261").Value = 100102005
$ws.Range("J261").Value = "Naranja"
$ws.Range("K261").Value = "Fukumoto"
$ws.Range("L261").Value = "Primera"
$ws.Range("M261").Value = 300
$ws.Range("N261").Value = 6000
$ws.Range("O261").Value = 6000
$ws.Range("P261").Value = 6000
$ws.Range("Q261").Value = "$/bandeja 15 kilos granel"
$ws.Range("R261").Value = "Provincia de Melipilla"
$ws.Range("S261").Value = 400
$ws.Range("T261").Value = 15

$ws.Range("A262").Value = 5
$ws.Range("B262").Value = "Macroferia Regional de Talca"
$ws.Range("C262").Value = "Maule"
$ws.Range("D262").Value = 44400
$ws.Range("E262").Value = 7
$ws.Range("F262").Value = "Fruta"
$ws.Range("G262").Value = 100102
$ws.Range("H262").Value = "Cítricos"
$ws.Range("I262").Value = 100102005
$ws.Range("J262").Value = "Naranja"
$ws.Range("K262").Value = "Navel Late"
$ws.Range("L262").Value = "Primera"
$ws.Range("M262").Value = 320
$ws.Range("N262").Value = 6000
$ws.Range("O262").Value = 6000
$ws.Range("P262").Value = 6000
$ws.Range("Q262").Value = "$/bandeja 15 kilos granel"
$ws.Range("R262").Value = "Región de O'Higgins"
$ws.Range("S262").Value = 400
$ws.Range("T262").Value = 15

# Make sure the date columns keep the workbook's existing date/time
# number format (style index 2 in styles.xml) for consistency with the
# rest of column D.
$dateFormat = $ws.Range("D254").NumberFormat
$ws.Range("D255").NumberFormat = $dateFormat
$ws.Range("D256").NumberFormat = $dateFormat
$ws.Range("D257").NumberFormat = $dateFormat
$ws.Range("D258").NumberFormat = $dateFormat
$ws.Range("D259").NumberFormat = $dateFormat
$ws.Range("D260").NumberFormat = $dateFormat
$ws.Range("D261").NumberFormat = $dateFormat
$ws.Range("D262").NumberFormat = $dateFormat
